# Update the "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3104
    5  = 2719
    9  = 1441
    11 = 63
    12 = 18
    13 = 1219
    14 = 6
    22 = 2643
    24 = 307
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
